$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.383.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.07%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.249.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.73%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.27%  "

# Row 6
$ws.Range("E6").Value = "  +0.59%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "64.09"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.65%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("E9").Value = "  -2.19%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0948"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.72%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.08%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.75%  "

# Row 13
$ws.Range("E13").Value = "  -1.03%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.585.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.70%  "

# Row 15
$ws.Range("E15").Value = "  -5.36%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.28%  "

# Row 17
$ws.Range("E17").Value = "  -2.16%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.249.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.65%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.251.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.82%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0965"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.47%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.98%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.30%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "246.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.06%  "

# Row 24
$ws.Range("E24").Value = "  +17.06%  "

# Row 25
$ws.Range("E25").Value = "  -0.06%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.85%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.27%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "173.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.71%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.80%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.35%  "

# Row 31
$ws.Range("E31").Value = "  +3.37%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.129"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.64%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.124"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.20%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.28%  "

# Row 35
$ws.Range("E35").Value = "  -0.82%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.43%  "

# Row 37
$ws.Range("E37").Value = "  -5.13%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.29%  "

# Row 39
$ws.Range("E39").Value = "  -2.07%  "

# Row 40
$ws.Range("E40").Value = "  -3.68%  "

# Row 41
$ws.Range("E41").Value = "  -0.08%  "

# Row 42
$ws.Range("E42").Value = "  +5.89%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.49%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.26%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.26%  "

# Row 46
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0937"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.79%  "

# Row 47
$ws.Range("E47").Value = "  -1.08%  "

# Row 48
$ws.Range("B48").Value = "Celestia"
$ws.Range("C48").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.12%  "

# Row 49
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.422.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.34%  "

# Row 50
$ws.Range("B50").Value = "TerraClassic"
$ws.Range("C50").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000204"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.61%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.18%  "
